$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrCK = New-Object "object[,]" 24,9
$arrCK[0,0] = 4.201225230854027
$arrCK[0,1] = 8.545511156068597
$arrCK[0,2] = 13.40546389943834
$arrCK[0,3] = 35.79040185421714
$arrCK[0,4] = 39.81654225435312
$arrCK[0,5] = 16.83376209536237
$arrCK[0,6] = 24.9037067305201
$arrCK[0,7] = 10.01716297997086
$arrCK[0,8] = 16.97629248222802
$arrCK[1,0] = 4.15200703316514
$arrCK[1,1] = 8.507928704994264
$arrCK[1,2] = 13.35495372795927
$arrCK[1,3] = 35.76448812944
$arrCK[1,4] = 39.7164337123366
$arrCK[1,5] = 16.87651492728228
$arrCK[1,6] = 24.91460065437727
$arrCK[1,7] = 10.02502440056451
$arrCK[1,8] = 16.52405358118132
$arrCK[2,0] = 4.121093246493413
$arrCK[2,1] = 8.485964049680858
$arrCK[2,2] = 13.32662578939525
$arrCK[2,3] = 35.75944728162774
$arrCK[2,4] = 39.6701548693966
$arrCK[2,5] = 16.9066182981614
$arrCK[2,6] = 24.92864064448177
$arrCK[2,7] = 10.03158554830292
$arrCK[2,8] = 16.24302889848279
$arrCK[3,0] = 4.108327747038018
$arrCK[3,1] = 8.477298703556974
$arrCK[3,2] = 13.31576518697892
$arrCK[3,3] = 35.76012386043017
$arrCK[3,4] = 39.65511901335965
$arrCK[3,5] = 16.91985134810363
$arrCK[3,6] = 24.93620541932035
$arrCK[3,7] = 10.03469522265723
$arrCK[3,8] = 16.12784170644222
$arrCK[4,0] = 4.10619809966227
$arrCK[4,1] = 8.475877238510545
$arrCK[4,2] = 13.31400327670054
$arrCK[4,3] = 35.76040101273685
$arrCK[4,4] = 39.65285323841935
$arrCK[4,5] = 16.92210691481243
$arrCK[4,6] = 24.93757270608898
$arrCK[4,7] = 10.03523790617007
$arrCK[4,8] = 16.10867993282672
$arrCK[5,0] = 4.120921757522122
$arrCK[5,1] = 8.485846022367676
$arrCK[5,2] = 13.32647654307138
$arrCK[5,3] = 35.75944535484288
$arrCK[5,4] = 39.66993660976883
$arrCK[5,5] = 16.90679285806311
$arrCK[5,6] = 24.92873520969909
$arrCK[5,7] = 10.03162572156025
$arrCK[5,8] = 16.24147790852281
$arrCK[6,0] = 4.184402628605886
$arrCK[6,1] = 8.532325823689813
$arrCK[6,2] = 13.38749563298324
$arrCK[6,3] = 35.77920978812396
$arrCK[6,4] = 39.77887652930372
$arrCK[6,5] = 16.84770172603917
$arrCK[6,6] = 24.9059345110788
$arrCK[6,7] = 10.01951362094821
$arrCK[6,8] = 16.82115727867151
$arrCK[7,0] = 4.303106848228104
$arrCK[7,1] = 8.632010842602091
$arrCK[7,2] = 13.52807704337869
$arrCK[7,3] = 35.90423161261916
$arrCK[7,4] = 40.11259107097963
$arrCK[7,5] = 16.76253186042168
$arrCK[7,6] = 24.91972666868882
$arrCK[7,7] = 10.00952391638546
$arrCK[7,8] = 17.92416010350993
$arrCK[8,0] = 4.386435575558495
$arrCK[8,1] = 8.710064414146116
$arrCK[8,2] = 13.64355290033086
$arrCK[8,3] = 36.04853848307786
$arrCK[8,4] = 40.4300381341503
$arrCK[8,5] = 16.71885509077233
$arrCK[8,6] = 24.96568989212571
$arrCK[8,7] = 10.01057143606198
$arrCK[8,8] = 18.70508892460769
$arrCK[9,0] = 4.423428487765154
$arrCK[9,1] = 8.746530036647281
$arrCK[9,2] = 13.69859789138967
$arrCK[9,3] = 36.12549584017076
$arrCK[9,4] = 40.58983168269582
$arrCK[9,5] = 16.70312306119138
$arrCK[9,6] = 24.99438794211925
$arrCK[9,7] = 10.0128659347099
$arrCK[9,8] = 19.05230379313959
$arrCK[10,0] = 4.437299512059589
$arrCK[10,1] = 8.760468627984999
$arrCK[10,2] = 13.71979114225464
$arrCK[10,3] = 36.1562530765403
$arrCK[10,4] = 40.65251923353205
$arrCK[10,5] = 16.69776307044217
$arrCK[10,6] = 25.00637353728267
$arrCK[10,7] = 10.01399568421534
$arrCK[10,8] = 19.18250703848173
$arrCK[11,0] = 4.434318342219163
$arrCK[11,1] = 8.757461057209721
$arrCK[11,2] = 13.71521148659309
$arrCK[11,3] = 36.14955732947228
$arrCK[11,4] = 40.63892212891909
$arrCK[11,5] = 16.69889083348473
$arrCK[11,6] = 25.0037425331371
$arrCK[11,7] = 10.01374078034175
$arrCK[11,8] = 19.15452420638607
$arrCK[12,0] = 4.424572452285163
$arrCK[12,1] = 8.747674217043343
$arrCK[12,2] = 13.70033455488957
$arrCK[12,3] = 36.12799396076068
$arrCK[12,4] = 40.59494560998409
$arrCK[12,5] = 16.70267010530539
$arrCK[12,6] = 24.99535160777497
$arrCK[12,7] = 10.01295365526295
$arrCK[12,8] = 19.06304200772728
$arrCK[13,0] = 4.41858474719351
$arrCK[13,1] = 8.741696164567744
$arrCK[13,2] = 13.69126705680012
$arrCK[13,3] = 36.11499573378103
$arrCK[13,4] = 40.56829112941983
$arrCK[13,5] = 16.70506288713048
$arrCK[13,6] = 24.99035747785124
$arrCK[13,7] = 10.01250547298515
$arrCK[13,8] = 19.0068363937107
$arrCK[14,0] = 4.383999097437654
$arrCK[14,1] = 8.707699911046044
$arrCK[14,2] = 13.64000510851579
$arrCK[14,3] = 36.04373583586818
$arrCK[14,4] = 40.41990194781695
$arrCK[14,5] = 16.71996667900321
$arrCK[14,6] = 24.96397102207166
$arrCK[14,7] = 10.01045802374051
$arrCK[14,8] = 18.68222508268621
$arrCK[15,0] = 4.362543554824104
$arrCK[15,1] = 8.68708421612901
$arrCK[15,2] = 13.609192751137
$arrCK[15,3] = 36.00291001093476
$arrCK[15,4] = 40.3327870957072
$arrCK[15,5] = 16.73017106167764
$arrCK[15,6] = 24.94977792984191
$arrCK[15,7] = 10.00966726110249
$arrCK[15,8] = 18.48093699920488
$arrCK[16,0] = 4.350117232857746
$arrCK[16,1] = 8.675317339683501
$arrCK[16,2] = 13.59170776712551
$arrCK[16,3] = 35.98049360830476
$arrCK[16,4] = 40.28413150576246
$arrCK[16,5] = 16.73642958598185
$arrCK[16,6] = 24.94234778617533
$arrCK[16,7] = 10.00938359388866
$arrCK[16,8] = 18.36441197269689
$arrCK[17,0] = 4.345895367615619
$arrCK[17,1] = 8.671349100377169
$arrCK[17,2] = 13.58582881366454
$arrCK[17,3] = 35.97308711184962
$arrCK[17,4] = 40.26790768532214
$arrCK[17,5] = 16.73861538531108
$arrCK[17,6] = 24.93995805255653
$arrCK[17,7] = 10.00931696010453
$arrCK[17,8] = 18.32483385255426
$arrCK[18,0] = 4.364836443714172
$arrCK[18,1] = 8.689269459598258
$arrCK[18,2] = 13.61244829193015
$arrCK[18,3] = 36.0071457860918
$arrCK[18,4] = 40.3419107314143
$arrCK[18,5] = 16.72904448468101
$arrCK[18,6] = 24.95121291891328
$arrCK[18,7] = 10.00973372843301
$arrCK[18,8] = 18.50244295962579
$arrCK[19,0] = 4.427438831737624
$arrCK[19,1] = 8.750545389826794
$arrCK[19,2] = 13.70469490493821
$arrCK[19,3] = 36.13428391348445
$arrCK[19,4] = 40.60780379588884
$arrCK[19,5] = 16.70154380958911
$arrCK[19,6] = 24.99778589775373
$arrCK[19,7] = 10.01317777806098
$arrCK[19,8] = 19.089948207355
$arrCK[20,0] = 4.467549685406738
$arrCK[20,1] = 8.791345785948621
$arrCK[20,2] = 13.76701082453619
$arrCK[20,3] = 36.22678332595341
$arrCK[20,4] = 40.79425092979115
$arrCK[20,5] = 16.68705357431622
$arrCK[20,6] = 25.03474033134549
$arrCK[20,7] = 10.01694887310087
$arrCK[20,8] = 19.46640912255597
$arrCK[21,0] = 4.446217202816676
$arrCK[21,1] = 8.76950363571131
$arrCK[21,2] = 13.73357046803448
$arrCK[21,3] = 36.17655830137576
$arrCK[21,4] = 40.6935941094553
$arrCK[21,5] = 16.69446780931273
$arrCK[21,6] = 25.0144217662348
$arrCK[21,7] = 10.01479728806567
$arrCK[21,8] = 19.26621000694531
$arrCK[22,0] = 4.363800112394537
$arrCK[22,1] = 8.688281244542766
$arrCK[22,2] = 13.61097574629388
$arrCK[22,3] = 36.00522750537249
$arrCK[22,4] = 40.33778148535021
$arrCK[22,5] = 16.72955258962302
$arrCK[22,6] = 24.9505618874287
$arrCK[22,7] = 10.00970314602241
$arrCK[22,8] = 18.4927226028567
$arrCK[23,0] = 4.271650901383677
$arrCK[23,1] = 8.60416850301149
$arrCK[23,2] = 13.48786144844369
$arrCK[23,3] = 35.86117871610205
$arrCK[23,4] = 40.00953203336389
$arrCK[23,5] = 16.78226605480069
$arrCK[23,6] = 24.90971010156083
$arrCK[23,7] = 10.0107528286316
$arrCK[23,8] = 17.63034035192478
$ws.Range("C2:K25").Value = $arrCK

$arrN = New-Object "object[,]" 24,1
$arrN[0,0] = 17.66049602683462
$arrN[1,0] = 17.71828516742664
$arrN[2,0] = 17.75571013592222
$arrN[3,0] = 17.77145067466774
$arrN[4,0] = 17.77409398237507
$arrN[5,0] = 17.755920434571
$arrN[6,0] = 17.68001932349955
$arrN[7,0] = 17.54653491147026
$arrN[8,0] = 17.45775290856393
$arrN[9,0] = 17.4193655523462
$arrN[10,0] = 17.40511577247949
$arrN[11,0] = 17.40817198175454
$arrN[12,0] = 17.41818747407992
$arrN[13,0] = 17.42435955953888
$arrN[14,0] = 17.46030177859813
$arrN[15,0] = 17.48286275343443
$arrN[16,0] = 17.49602753418657
$arrN[17,0] = 17.50051727284736
$arrN[18,0] = 17.48044161717389
$arrN[19,0] = 17.41523790801328
$arrN[20,0] = 17.37429409295845
$arrN[21,0] = 17.39599401605613
$arrN[22,0] = 17.48153560795716
$arrN[23,0] = 17.58100955630061
$ws.Range("N2:N25").Value = $arrN
